$wb = $excel.ActiveWorkbook

# --- Rename the existing sheet ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Top AL and NL Data"

# --- Remove the now-separated "Year WAR Ranking" rows (45-55) from sheet 1 ---
$ws1.Range("A45:AF55").EntireRow.Delete()

# --- Add the new sheet right after the first one ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Year WAR Ranking"

# Header row (copy the same header style used on sheet 1)
$ws2.Range("A1").Value = "Year"
$ws2.Range("B1").Value = "WAR/pos"
$ws1.Range("A1").Copy()
$ws2.Range("A1:B1").PasteSpecial(-4122)

# Data rows (Year, WAR/pos)
$data = @(
    @(2011, 595.8),
    @(2012, 595.3),
    @(2009, 594.8),
    @(2010, 594),
    @(2007, 593.6),
    @(2008, 593),
    @(2015, 592.8),
    @(2013, 592),
    @(2016, 591.6),
    @(2014, 590.6)
)

$r = 2
foreach ($row in $data) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

$ws1.Activate()
